$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 108833
$ws.Range("E2").Value = 3198
$ws.Range("F2").Value = 3727
$ws.Range("G2").Value = 1517
$ws.Range("H2").Value = 867
$ws.Range("I2").Value = 450
$ws.Range("J2").Value = 417
$ws.Range("K2").Value = 103357
$ws.Range("L2").Value = 73026
$ws.Range("M2").Value = 30331
$ws.Range("N2").Value = 24051
$ws.Range("O2").Value = 6280
$ws.Range("P2").Value = 1610
$ws.Range("Q2").Value = 4092
$ws.Range("R2").Value = -2629
$ws.Range("S2").Value = -523
$ws.Range("T2").Value = 2292
$ws.Range("U2").Value = 1799
$ws.Range("V2").Value = 47807
$ws.Range("W2").Value = 2.94
$ws.Range("X2").Value = 0.8
$ws.Range("Y2").Value = 1.87
$ws.Range("Z2").Value = 0.85
$ws.Range("AA2").Value = 240.76
$ws.Range("AB2").Value = 1461.1
$ws.Range("AC2").Value = 1397
$ws.Range("AD2").Value = 38.86
$ws.Range("AE2").Value = 86616
$ws.Range("AF2").Value = 0.63
$ws.Range("AG2").Value = 1250
$ws.Range("AH2").Value = 2.3
$ws.Range("AI2").Value = 77.14
$ws.Range("AJ2").Value = 32200000

# Row 3
$ws.Range("D3").Value = 99997
$ws.Range("E3").Value = 3143
$ws.Range("F3").Value = 2720
$ws.Range("G3").Value = 170
$ws.Range("H3").Value = -735
$ws.Range("I3").Value = -987
$ws.Range("J3").Value = 252
$ws.Range("K3").Value = 98815
$ws.Range("L3").Value = 69279
$ws.Range("M3").Value = 29536
$ws.Range("N3").Value = 23026
$ws.Range("O3").Value = 6510
$ws.Range("P3").Value = 1610
$ws.Range("Q3").Value = 4641
$ws.Range("R3").Value = -3428
$ws.Range("S3").Value = -1393
$ws.Range("T3").Value = 2535
$ws.Range("U3").Value = 2106
$ws.Range("V3").Value = 47030
$ws.Range("W3").Value = 3.14
$ws.Range("X3").Value = -0.74
$ws.Range("Y3").Value = -4.19
$ws.Range("Z3").Value = -0.73
$ws.Range("AA3").Value = 234.56
$ws.Range("AB3").Value = 1380.67
$ws.Range("AC3").Value = -3064
$ws.Range("AD3").Value = -13.17
$ws.Range("AE3").Value = 82926
$ws.Range("AF3").Value = 0.49
$ws.Range("AG3").Value = 1250
$ws.Range("AH3").Value = 3.1
$ws.Range("AI3").Value = -35.18
$ws.Range("AJ3").Value = 32200000

# Row 4
$ws.Range("D4").Value = 85075
$ws.Range("E4").Value = 2793
$ws.Range("F4").Value = 4592
$ws.Range("G4").Value = 2030
$ws.Range("H4").Value = 2151
$ws.Range("I4").Value = 1682
$ws.Range("J4").Value = 470
$ws.Range("K4").Value = 97130
$ws.Range("L4").Value = 64578
$ws.Range("M4").Value = 32552
$ws.Range("N4").Value = 24515
$ws.Range("O4").Value = 8038
$ws.Range("P4").Value = 1610
$ws.Range("Q4").Value = 8079
$ws.Range("R4").Value = -2359
$ws.Range("S4").Value = -6119
$ws.Range("T4").Value = 1671
$ws.Range("U4").Value = 6409
$ws.Range("V4").Value = 40471
$ws.Range("W4").Value = 3.28
$ws.Range("X4").Value = 2.53
$ws.Range("Y4").Value = 7.07
$ws.Range("Z4").Value = 2.2
$ws.Range("AA4").Value = 198.38
$ws.Range("AB4").Value = 1474.75
$ws.Range("AC4").Value = 5223
$ws.Range("AD4").Value = 11.35
$ws.Range("AE4").Value = 88286
$ws.Range("AF4").Value = 0.67
$ws.Range("AG4").Value = 1250
$ws.Range("AH4").Value = 2.11
$ws.Range("AI4").Value = 20.64
$ws.Range("AJ4").Value = 32200000

# Row 5
$ws.Range("D5").Value = 94183
$ws.Range("E5").Value = 4025
$ws.Range("F5").Value = 5274
$ws.Range("G5").Value = 3742
$ws.Range("H5").Value = 3648
$ws.Range("I5").Value = 2969
$ws.Range("J5").Value = 680
$ws.Range("K5").Value = 99789
$ws.Range("L5").Value = 63468
$ws.Range("M5").Value = 36321
$ws.Range("N5").Value = 27217
$ws.Range("O5").Value = 9104
$ws.Range("P5").Value = 1610
$ws.Range("Q5").Value = 2374
$ws.Range("R5").Value = -4030
$ws.Range("S5").Value = 1490
$ws.Range("T5").Value = 1793
$ws.Range("U5").Value = 582
$ws.Range("V5").Value = 37908
$ws.Range("W5").Value = 4.27
$ws.Range("X5").Value = 3.87
$ws.Range("Y5").Value = 11.48
$ws.Range("Z5").Value = 3.71
$ws.Range("AA5").Value = 174.74
$ws.Range("AB5").Value = 1662.86
$ws.Range("AC5").Value = 9219
$ws.Range("AD5").Value = 7.86
$ws.Range("AE5").Value = 98018
$ws.Range("AF5").Value = 0.74
$ws.Range("AG5").Value = 1250
$ws.Range("AH5").Value = 1.72
$ws.Range("AI5").Value = 11.69
$ws.Range("AJ5").Value = 32200000

# Row 6
$ws.Range("D6").Value = 101102
$ws.Range("E6").Value = 3543
$ws.Range("F6").Value = 5091
$ws.Range("G6").Value = 2751
$ws.Range("H6").Value = 4877
$ws.Range("I6").Value = 4041
$ws.Range("K6").Value = 104903
$ws.Range("L6").Value = 61307
$ws.Range("M6").Value = 43596
$ws.Range("N6").Value = 31636
$ws.Range("P6").Value = 1610
$ws.Range("Q6").Value = 499
$ws.Range("R6").Value = 3667
$ws.Range("S6").Value = -2315
$ws.Range("T6").Value = 2985
$ws.Range("U6").Value = -2486
$ws.Range("V6").Value = 38442
$ws.Range("W6").Value = 3.5
$ws.Range("X6").Value = 4.82
$ws.Range("Y6").Value = 13.73
$ws.Range("Z6").Value = 4.76
$ws.Range("AA6").Value = 140.63
$ws.Range("AB6").Value = 1942.82
$ws.Range("AC6").Value = 12551
$ws.Range("AD6").Value = 3.91
$ws.Range("AE6").Value = 113934
$ws.Range("AF6").Value = 0.43
$ws.Range("AI6").Value = 8.59
$ws.Range("AJ6").Value = 32200000

# Row 7
$ws.Range("D7").Value = 100946
$ws.Range("E7").Value = 3843
$ws.Range("G7").Value = 2751
$ws.Range("H7").Value = 1873
$ws.Range("I7").Value = 1326
$ws.Range("K7").Value = 109124
$ws.Range("L7").Value = 64729
$ws.Range("M7").Value = 44395
$ws.Range("N7").Value = 32587
$ws.Range("P7").Value = 1610
$ws.Range("Q7").Value = 4583
$ws.Range("R7").Value = -1338
$ws.Range("S7").Value = -555
$ws.Range("T7").Value = 2609
$ws.Range("U7").Value = 1786
$ws.Range("W7").Value = 3.81
$ws.Range("X7").Value = 1.86
$ws.Range("Y7").Value = 4.13
$ws.Range("Z7").Value = 1.75
$ws.Range("AA7").Value = 145.8
$ws.Range("AC7").Value = 4118
$ws.Range("AD7").Value = 9.68
$ws.Range("AE7").Value = 117357
$ws.Range("AF7").Value = 0.34
$ws.Range("AG7").Value = 1250
$ws.Range("AH7").Value = 3.14
$ws.Range("AI7").Value = 30.35

# Row 8
$ws.Range("D8").Value = 103686
$ws.Range("E8").Value = 4561
$ws.Range("G8").Value = 3551
$ws.Range("H8").Value = 2614
$ws.Range("I8").Value = 1818
$ws.Range("K8").Value = 111252
$ws.Range("L8").Value = 64668
$ws.Range("M8").Value = 46584
$ws.Range("N8").Value = 34013
$ws.Range("P8").Value = 1610
$ws.Range("Q8").Value = 3998
$ws.Range("R8").Value = -1530
$ws.Range("S8").Value = -1365
$ws.Range("T8").Value = 1300
$ws.Range("U8").Value = 4122
$ws.Range("W8").Value = 4.4
$ws.Range("X8").Value = 2.52
$ws.Range("Y8").Value = 5.46
$ws.Range("Z8").Value = 2.37
$ws.Range("AA8").Value = 138.82
$ws.Range("AC8").Value = 5647
$ws.Range("AD8").Value = 7.06
$ws.Range("AE8").Value = 122494
$ws.Range("AF8").Value = 0.33
$ws.Range("AG8").Value = 1250
$ws.Range("AH8").Value = 3.14
$ws.Range("AI8").Value = 22.14

# Row 9
$ws.Range("D9").Value = 106407
$ws.Range("E9").Value = 4697
$ws.Range("G9").Value = 3918
$ws.Range("H9").Value = 2804
$ws.Range("I9").Value = 2003
$ws.Range("K9").Value = 113912
$ws.Range("L9").Value = 64800
$ws.Range("M9").Value = 49112
$ws.Range("N9").Value = 35742
$ws.Range("P9").Value = 1610
$ws.Range("Q9").Value = 4632
$ws.Range("R9").Value = -1705
$ws.Range("S9").Value = -1391
$ws.Range("T9").Value = 1275
$ws.Range("U9").Value = 3155
$ws.Range("W9").Value = 4.41
$ws.Range("X9").Value = 2.63
$ws.Range("Y9").Value = 5.74
$ws.Range("Z9").Value = 2.49
$ws.Range("AA9").Value = 131.94
$ws.Range("AC9").Value = 6219
$ws.Range("AD9").Value = 6.41
$ws.Range("AE9").Value = 128722
$ws.Range("AF9").Value = 0.31
$ws.Range("AG9").Value = 1250
$ws.Range("AH9").Value = 3.14
$ws.Range("AI9").Value = 20.1

# Remove AG6 and AH6 entirely (cells deleted in target)
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
